$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late"), shifting the
# existing N/O/P data right by one column (N->O, O->P, P->Q).
$ws.Columns("N").Insert()

# Match the new column's width to column M's width, the way Excel copies
# the width of the preceding column when a column is inserted.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and restore the
# remembered selection on it.
$ws.Activate()
$ws.Range("J22").Select()
